$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header-row shared strings:
#    "<name>_old" -> "<name>_FV2310"
#    "<name>_new" -> "<name>_FV2404"
$renameMap = @{
    "Segmentname_old"        = "Segmentname_FV2310";
    "Segmentgruppe_old"      = "Segmentgruppe_FV2310";
    "Segment_old"            = "Segment_FV2310";
    "Datenelement_old"       = "Datenelement_FV2310";
    "Segment ID_old"         = "Segment ID_FV2310";
    "Code_old"               = "Code_FV2310";
    "Qualifier_old"          = "Qualifier_FV2310";
    "Beschreibung_old"       = "Beschreibung_FV2310";
    "Bedingungsausdruck_old" = "Bedingungsausdruck_FV2310";
    "Bedingung_old"          = "Bedingung_FV2310";
    "Segmentname_new"        = "Segmentname_FV2404";
    "Segmentgruppe_new"      = "Segmentgruppe_FV2404";
    "Segment_new"            = "Segment_FV2404";
    "Datenelement_new"       = "Datenelement_FV2404";
    "Segment ID_new"         = "Segment ID_FV2404";
    "Code_new"               = "Code_FV2404";
    "Qualifier_new"          = "Qualifier_FV2404";
    "Beschreibung_new"       = "Beschreibung_FV2404";
    "Bedingungsausdruck_new" = "Bedingungsausdruck_FV2404";
    "Bedingung_new"          = "Bedingung_FV2404";
}

$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($renameMap.ContainsKey($val)) {
        $cell.Value = $renameMap[$val]
    }
}

# 2) Turn the populated range into an Excel Table (ListObject) named "Table1"
$lastRow = $ws.UsedRange.Rows.Count
$tblRange = $ws.Range("A1:U$lastRow")
$tbl = $ws.ListObjects.Add(1, $tblRange, $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split/freeze pane below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
